$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("pages_with_translation")

# "Spanish" -> "Español" wherever it appears (column E, language toggle values)
$usedRange = $ws.UsedRange
for ($r = 1; $r -le $usedRange.Rows.Count; $r++) {
    $cell = $ws.Cells.Item($r, 5)
    if ($cell.Value() -eq "Spanish") {
        $cell.Value = "Español"
    }
}

# Bold the header row
$ws.Range("A1:E1").Font.Bold = $true

# Widen column E
$ws.Columns.Item(5).ColumnWidth = 17.67

# Move active selection to E9
$ws.Range("E9").Select()
